$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: update title and link
$ws.Range("D32").Value = "RESTful API 파이썬 구축 예시 (feat. Flask)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/416"

# Row 44: update title and link
$ws.Range("D44").Value = "Microsoft 365 Copilot 정리"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/126"

# Row 51: update title and link
$ws.Range("D51").Value = "[python] 얕은 복사와 깊은 복사, 이거 모르면 큰일남"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EC%96%95%EC%9D%80-%EB%B3%B5%EC%82%AC%EC%99%80-%EA%B9%8A%EC%9D%80-%EB%B3%B5%EC%82%AC-%EC%9D%B4%EA%B1%B0-%EB%AA%A8%EB%A5%B4%EB%A9%B4-%ED%81%B0%EC%9D%BC%EB%82%A8"
